# Applies the dated-worksheet update: new date and new multiplication problems.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-01-05 Sunday"; New = "2025-01-06 Monday" },
    @{ Old = "23×51="; New = "80×50=" },
    @{ Old = "24×13="; New = "88×63=" },
    @{ Old = "85×74="; New = "73×49=" },
    @{ Old = "75×63="; New = "60×78=" },
    @{ Old = "69×66="; New = "30×92=" },
    @{ Old = "18×82="; New = "45×78=" },
    @{ Old = "59×95="; New = "25×76=" },
    @{ Old = "21×11="; New = "13×42=" },
    @{ Old = "32×90="; New = "30×60=" },
    @{ Old = "90×61="; New = "17×27=" },
    @{ Old = "74×49="; New = "58×85=" },
    @{ Old = "11×75="; New = "47×27=" },
    @{ Old = "62×70="; New = "87×59=" },
    @{ Old = "29×95="; New = "97×58=" },
    @{ Old = "20×42="; New = "16×73=" },
    @{ Old = "59×13="; New = "62×88=" },
    @{ Old = "37×81="; New = "30×19=" },
    @{ Old = "88×17="; New = "52×91=" },
    @{ Old = "62×76="; New = "35×46=" },
    @{ Old = "70×56="; New = "95×87=" },
    @{ Old = "58×58="; New = "40×70=" },
    @{ Old = "43×51="; New = "29×26=" },
    @{ Old = "57×19="; New = "27×79=" },
    @{ Old = "60×33="; New = "86×50=" },
    @{ Old = "78×36="; New = "91×19=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
